$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# "Volume 31   Number  12" -> "Volume 31   Number  13"
$ws.Range("A8").Value = "Volume 31   Number  13"
# "Report Covering the Week  3/18/2024  Through  3/24/2024" -> 3/25/2024 .. 3/31/2024
$ws.Range("C9").Value = "Report Covering the Week  3/25/2024  Through  3/31/2024"

# --- Crime statistics table updates (rows 15-28) ---
# Cells that change from a text placeholder ("0" / "***.*") to a real number, or
# vice-versa, need their number format copied from a same-row donor cell that
# already carries the target style, otherwise Excel would invent a brand new
# style entry. PasteSpecial(xlPasteFormats = -4122) copies formatting only.

# C15: text "0"
$ws.Range("C15").Value = "'0"
# D15: number 1 (style fix from donor F15)
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
# E15: number -100 (style fix from donor L15)
$ws.Range("L15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
# G15: number 1 (style fix from donor F15)
$ws.Range("F15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
# H15: number 0 (style fix from donor L15)
$ws.Range("L15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = -75
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -66.666666666666

# C16: text "0" (style fix from donor C15)
$ws.Range("C16").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 31
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = -29.545454545454
$ws.Range("L16").Value = -46.551724137931
$ws.Range("M16").Value = -16.216216216216
$ws.Range("N16").Value = -85.238095238095

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -66.666666666666
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = -46.341463414634
$ws.Range("L17").Value = -31.25
$ws.Range("M17").Value = 4.761904761904
$ws.Range("N17").Value = -54.166666666666

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -36
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = -27.631578947368
$ws.Range("L18").Value = -35.294117647058
$ws.Range("M18").Value = -3.508771929824
$ws.Range("N18").Value = -69.101123595505

$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -18.75
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = -16.279069767441
$ws.Range("I19").Value = 238
$ws.Range("J19").Value = 306
$ws.Range("K19").Value = -22.222222222222
$ws.Range("L19").Value = -3.643724696356
$ws.Range("M19").Value = 1.276595744680
$ws.Range("N19").Value = -56.25

$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = 71.428571428571
$ws.Range("L20").Value = 20
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -94.029850746268

$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -41.935483870967
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = -19.708029197080
$ws.Range("I21").Value = 359
$ws.Range("J21").Value = 476
$ws.Range("K21").Value = -24.579831932773
$ws.Range("L21").Value = -17.660550458715
$ws.Range("M21").Value = 0.279329608938
$ws.Range("N21").Value = -69.679054054054

# C22: number 1 (style fix from donor F22)
$ws.Range("F22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 13
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 8.333333333333
$ws.Range("M22").Value = -35

$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 39
$ws.Range("F24").Value = 131
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = 14.912280701754
$ws.Range("I24").Value = 409
$ws.Range("J24").Value = 419
$ws.Range("K24").Value = -2.386634844868
$ws.Range("L24").Value = 2.763819095477
$ws.Range("M24").Value = 17.528735632183

$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = 11.538461538461
$ws.Range("F25").Value = 109
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = 36.25
$ws.Range("I25").Value = 335
$ws.Range("J25").Value = 307
$ws.Range("K25").Value = 9.120521172638
$ws.Range("L25").Value = 0.299401197604

$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 87
$ws.Range("J26").Value = 99
$ws.Range("K26").Value = -12.121212121212
$ws.Range("L26").Value = 17.567567567567
$ws.Range("M26").Value = 77.551020408163

# D27: number 1 (style fix from donor F27)
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
# E27: number -100 (style fix from donor K27)
$ws.Range("K27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
# G27: number 1 (style fix from donor F27)
$ws.Range("F27").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
# H27: number 0 (style fix from donor K27)
$ws.Range("K27").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = -50

$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = -5.263157894736
$ws.Range("L28").Value = 12.5
